$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price / volume(1h) values.
# Numeric-looking values are forced to remain text (matching the
# original inlineStr cell type) by temporarily using a text number
# format, then resetting the style so no stray style id is left on
# the cell.

$ws.Range("D2").Value = '68.621.35'
$ws.Range("E2").Value = '  +2.02%  '
$ws.Range("D3").Value = '2.646.59'
$ws.Range("E3").Value = '  +1.35%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.23'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.51%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.84'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.84%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.546'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.13%  '
$ws.Range("D9").Value = '2.641.73'
$ws.Range("E9").Value = '  +1.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.137'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +12.14%  '
$ws.Range("E11").Value = '  -0.43%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.23'
$ws.Range("D12").Style = "Normal"
$ws.Range("E13").Value = '  +0.36%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.74'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000189'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.69%  '
$ws.Range("D16").Value = '3.120.46'
$ws.Range("E16").Value = '  +1.28%  '
$ws.Range("D17").Value = '68.454.86'
$ws.Range("E17").Value = '  +2.02%  '
$ws.Range("D18").Value = '2.644.59'
$ws.Range("E18").Value = '  +1.46%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.41'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.29%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '367.49'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.49'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.44%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.26'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.77%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.86'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.41%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.10'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.69%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '73.16'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +10.54%  '
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.85'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.57%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000105'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.74%  '
$ws.Range("E30").Value = '  +0.08%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '575.66'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.72%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.97'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.51%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.41'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.77%  '
$ws.Range("E34").Value = '  +2.70%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.997'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.22%  '
$ws.Range("E36").Value = '  +4.19%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.52'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.43%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '160.03'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.58%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.22'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.71%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.90'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.10%  '
$ws.Range("E41").Value = '  +0.14%  '
$ws.Range("E42").Value = '  +1.68%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.66'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.92%  '
$ws.Range("E44").Value = '  +4.98%  '
$ws.Range("D45").Value = '0.0₆0325'
$ws.Range("E45").Value = '  +10.63%  '
$ws.Range("E46").Value = '  +0.05%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '40.50'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.80%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '155.97'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.57%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.73'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.23%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.00'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.62%  '
$ws.Range("E51").Value = '  +0.26%  '
